# Update the cryptos worksheet with freshly "scraped" values.
# NOTE: the Price column (D) frequently holds values that *look* numeric
# (e.g. "142.85", "0.608"). Excel's smart-entry would silently convert a
# plain Range.Value assignment of such a string into a real number (and
# lose the original inline/shared-string text representation). To keep
# these as text - exactly like the source data - we briefly force the
# cell to Text format, assign the string, and then clear the formatting
# again so the cell's style is left untouched (matching the original
# default, un-styled cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "69.160.15"
$ws.Range("E2").Value = "  +1.26%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.893.07"
$ws.Range("E3").Value = "  -0.72%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.21%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +7.98%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "142.85"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.608"
$ws.Range("E7").Value = "  -2.48%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.717"
$ws.Range("E9").Value = "  -2.41%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.169"
$ws.Range("E10").Value = "  +0.37%  "

# Row 11 - ShibaInu
Set-TextValue $ws.Range("D11") "0.0000328"
$ws.Range("E11").Value = "  -6.17%  "

# Row 12 - Avalanche
Set-TextValue $ws.Range("D12") "41.81"
$ws.Range("E12").Value = "  -3.26%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "4.526.75"

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "10.15"
$ws.Range("E14").Value = "  -5.40%  "

# Row 15 - WrappedEther
Set-TextValue $ws.Range("D15") "3.917.81"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  -0.57%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +6.23%  "

# Row 18 - Uniswap
Set-TextValue $ws.Range("D18") "13.69"
$ws.Range("E18").Value = "  -4.27%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "19.59"
$ws.Range("E19").Value = "  -2.64%  "

# Row 20 - WrappedBTC
Set-TextValue $ws.Range("D20") "69.018.10"
$ws.Range("E20").Value = "  +0.99%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "424.44"
$ws.Range("E21").Value = "  -1.78%  "

# Row 22 - ImmutableX
Set-TextValue $ws.Range("D22") "3.31"
$ws.Range("E22").Value = "  -5.78%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D23") "14.10"
$ws.Range("E23").Value = "  -7.07%  "

# Row 24 - was Litecoin, now PancakeSwap
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D24") "4.08"
$ws.Range("E24").Value = "  +9.69%  "

# Row 25 - was PancakeSwap, now Litecoin
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "87.36"
$ws.Range("E25").Value = "  -1.37%  "

# Row 26 - RenderToken
Set-TextValue $ws.Range("D26") "11.41"
$ws.Range("E26").Value = "  -3.11%  "

# Row 27 - Filecoin
Set-TextValue $ws.Range("D27") "10.50"
$ws.Range("E27").Value = "  -6.36%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "35.75"
$ws.Range("E28").Value = "  -5.64%  "

# Row 29 - Bittensor
Set-TextValue $ws.Range("D29") "693.51"
$ws.Range("E29").Value = "  -3.52%  "

# Row 30 - Cosmos
Set-TextValue $ws.Range("D30") "13.02"
$ws.Range("E30").Value = "  -5.37%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -4.86%  "

# Row 32 - Toncoin
Set-TextValue $ws.Range("D32") "2.79"
$ws.Range("E32").Value = "  -4.39%  "

# Row 33 - OKB
Set-TextValue $ws.Range("D33") "68.49"
$ws.Range("E33").Value = "  +12.23%  "

# Row 34 - TheGraph
Set-TextValue $ws.Range("D34") "0.443"
$ws.Range("E34").Value = "  +11.46%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -4.80%  "

# Row 36 - InjectiveProtocol
Set-TextValue $ws.Range("D36") "39.78"
$ws.Range("E36").Value = "  -4.80%  "

# Row 37 - PEPE
$ws.Range("E37").Value = "  -9.15%  "

# Row 38 - Dai
Set-TextValue $ws.Range("D38") "0.999"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.08%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.02%  "

# Row 41 - VeChain
Set-TextValue $ws.Range("D41") "0.0476"
$ws.Range("E41").Value = "  -3.46%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +1.61%  "

# Row 43 - Fetch.AI
Set-TextValue $ws.Range("D43") "2.73"
$ws.Range("E43").Value = "  -9.32%  "

# Row 44 - ThetaToken
$ws.Range("E44").Value = "  -6.62%  "

# Row 45 - ApeXProtocol
$ws.Range("E45").Value = "  -0.18%  "

# Row 46 - Stellar
Set-TextValue $ws.Range("D46") "0.139"
$ws.Range("E46").Value = "  -2.26%  "

# Row 47 - Stacks
Set-TextValue $ws.Range("D47") "3.02"
$ws.Range("E47").Value = "  +7.30%  "

# Row 48 - LidoDAOToken
Set-TextValue $ws.Range("D48") "3.25"
$ws.Range("E48").Value = "  -5.05%  "

# Row 49 - Monero
Set-TextValue $ws.Range("D49") "142.40"
$ws.Range("E49").Value = "  -1.77%  "

# Row 50 - ARBITRUM
Set-TextValue $ws.Range("D50") "2.03"
$ws.Range("E50").Value = "  -5.32%  "

# Row 51 - was EnergySwap, now BabyDogeCoin
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.0₆0325"
$ws.Range("E51").Value = "  -4.60%  "
